# Composite vertical tailplane sizing - add SPAR V2 sheet, finalize D-CELL V2 / Panel V2
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Panel V2: the "M_max" check row (A16:C16) is being moved out to the new
#    SPAR V2 sheet - clear it here.
# ---------------------------------------------------------------------------
$panel = $wb.Worksheets.Item("Panel V2")
$panel.Range("A16:C16").ClearContents()
$panel.Range("B10").Select()

# ---------------------------------------------------------------------------
# 2. D-CELL V2: add the "Good!" labels next to the existing panel checks.
# ---------------------------------------------------------------------------
$dcell = $wb.Worksheets.Item("D-CELL V2")
$dcell.Range("D21").Value = "Min E"
$dcell.Range("D22").Value = "Good! "
$dcell.Range("D23").Value = "Good!"
$dcell.Range("B23").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "SPAR V2" worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$spar = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$spar.Name = "SPAR V2"

# ---- Moment / geometry block -----------------------------------------------
$spar.Range("A1").Value = "M_max"
$spar.Range("B1").Value = -600.22580000000005
$spar.Range("C1").Value = "Nm"
$spar.Range("D1").Value = "b max"
$spar.Range("E1").Value = 116.1
$spar.Range("G1").Value = "Facesheets"
$spar.Range("K1").Value = "Core (foam)"

$spar.Range("D2").Value = "b min [mm]"
$spar.Range("E2").Value = 81.8
$spar.Range("G2").Value = "E [Pa]"
$spar.Range("H2").Value = 72000000000
$spar.Range("H2").NumberFormat = "0.00E+00"
$spar.Range("K2").Value = "E [Pa]"
$spar.Range("L2").Value = 484000
$spar.Range("L2").NumberFormat = "0.00E+00"

$spar.Range("G3").Value = "t [mm]"
$spar.Range("H3").Value = 1
$spar.Range("K3").Value = "c max [mm]"
$spar.Range("L3").Formula = "=E1-H3"
$spar.Range("N3").Value = "c/t"
$spar.Range("O3").Formula = "=L3/H3"

$spar.Range("A4").Value = "L [m]"
$spar.Range("B4").Value = 1.1943999999999999
$spar.Range("G4").Value = "Al-7075, O"
$spar.Range("K4").Value = "c min [mm]"
$spar.Range("L4").Formula = "=E2-H3"
$spar.Range("O4").Formula = "=L4/H3"

$spar.Range("A5").Value = "b [m]"
$spar.Range("B5").Value = 0.1

$spar.Range("A6").Value = "c max"
$spar.Range("B6").Formula = "=L3*10^(-3)"
$spar.Range("C6").Value = "d max"
$spar.Range("D6").Formula = "=B6+B8"
$spar.Range("G6").Formula = "=B7/B4"

$spar.Range("A7").Value = "c min"
$spar.Range("B7").Formula = "=L4*10^(-3)"
$spar.Range("C7").Value = "d min"
$spar.Range("D7").Formula = "=B7+B8"
$spar.Range("G7").Formula = "=B8/B7"

$spar.Range("A8").Value = "t"
$spar.Range("B8").Formula = "=H3*10^(-3)"
$spar.Range("G8").Formula = "=B8/B6"

$spar.Range("A9").Value = "E_f"
$spar.Range("B9").Formula = "=H2"
$spar.Range("B9").NumberFormat = "0.00E+00"

$spar.Range("A10").Value = "E_c"
$spar.Range("B10").Formula = "=L2"
$spar.Range("B10").NumberFormat = "0.00E+00"
$spar.Range("E10").Value = "sigma_f"
$spar.Range("F10").Formula = "=ABS(B1)/(B5*B8*(B6+B8))"

$spar.Range("A12").Value = "sigma_ft"
$spar.Range("B12").Value = 214000000
$spar.Range("B12").NumberFormat = "0.00E+00"

$spar.Range("A13").Value = "sigma_fc"
$spar.Range("B13").Value = 100000000
$spar.Range("B13").NumberFormat = "0.00E+00"

$spar.Range("A14").Value = "tau_c"
$spar.Range("B14").Value = 50000
$spar.Range("B14").NumberFormat = "0.00E+00"
$spar.Range("C14").Value = "Shear strength"

$spar.Range("A15").Value = "sigma_core"
$spar.Range("B15").Value = 100000
$spar.Range("B15").NumberFormat = "0.00E+00"
$spar.Range("C15").Value = "Compressive strength of foam/ core"

# ---- Failure-load checks (red, highlighted cells) --------------------------
$spar.Range("A18").Font.Color = 255
$spar.Range("A18").Value = "P_FS_T"
$spar.Range("B18").Font.Color = 255
$spar.Range("B18").NumberFormat = "0.00E+00"
$spar.Range("B18").Formula = "=D7*4*B5*B8*B12/B4"
$spar.Range("D18").Value = "(EI)_eq"
$spar.Range("E18").Formula = "=B5*B8*B7^2*B9/2"
$spar.Range("E18").NumberFormat = "0.00E+00"

$spar.Range("A19").Font.Color = 255
$spar.Range("A19").Value = "P_FS_C"
$spar.Range("B19").Font.Color = 255
$spar.Range("B19").NumberFormat = "0.00E+00"
$spar.Range("B19").Formula = "=4*D7*B5*B8*B13/B4"
$spar.Range("D19").Value = "G_c"
$spar.Range("E19").Formula = "=350000"
$spar.Range("E19").NumberFormat = "0.00E+00"
$spar.Range("F19").Value = "Shear modulus of core"

$spar.Range("A20").Font.Color = 255
$spar.Range("B20").Font.Color = 255

$spar.Range("A21").Font.Color = 255
$spar.Range("A21").Value = "P_IN"
$spar.Range("B21").Font.Color = 255
$spar.Range("B21").NumberFormat = "0.00E+00"
$spar.Range("B21").Formula = "=B5*B8*((PI()^2*D7*B9*B15^2)/(3*B4))^(1/3)"
$spar.Range("D21").Value = "P_E"
$spar.Range("E21").Formula = "=4*PI()^2*E18/(B4^2)"
$spar.Range("E21").NumberFormat = "0.00E+00"
$spar.Range("F21").Value = "Euler buckling load"

$spar.Range("A22").Font.Color = 255
$spar.Range("B22").Font.Color = 255
$spar.Range("D22").Value = "P_s"
$spar.Range("E22").Formula = "=B5*B7*E19"
$spar.Range("E22").NumberFormat = "0.00E+00"
$spar.Range("F22").Value = "shear stiffness of core"

$spar.Range("A23").Font.Color = 255
$spar.Range("A23").Value = "P_CS"
$spar.Range("B23").Font.Color = 255
$spar.Range("B23").NumberFormat = "0.00E+00"
$spar.Range("B23").Formula = "=2*B5*D7*B14"

$spar.Range("A24").Font.Color = 255
$spar.Range("B24").Font.Color = 255
$spar.Range("D24").Value = "P_cr"
$spar.Range("E24").Formula = "=1/((1/E21)+(1/E22))"
$spar.Range("E24").NumberFormat = "0.00E+00"
$spar.Range("F24").Value = "Combined Collapse load"

$spar.Range("D27").Value = "sigma_fw"
$spar.Range("E27").Formula = "=0.5*(B9*B10*E19)^(1/3)"
$spar.Range("E27").NumberFormat = "0.00E+00"
$spar.Range("F27").Value = "Face wrinkling"

# ---- Cosmetics: column width, page setup, view ------------------------------
$spar.Columns.Item(6).ColumnWidth = 10.498697916666666
$spar.PageSetup.PaperSize = 9
$spar.PageSetup.Orientation = 1
$spar.Range("C32").Select()

Write-Host "Edit script completed"
